$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 315.54544
$ws.Range("I101").Value = 349.57144
$ws.Range("J101").Value = 256
$ws.Range("K101").Value = 1048.71432
$ws.Range("L101").Value = 768
$ws.Range("M101").Value = 573.28568
$ws.Range("N101").Value = -4012
$ws.Range("H135").Value = 792.9474
$ws.Range("I135").Value = 592.1177
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 5329.0593
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -2794.0593
$ws.Range("N135").Value = -27570
$ws.Range("H137").Value = 2142.647
$ws.Range("I137").Value = 2120.3125
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 6360.9375
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -3810.9375
$ws.Range("N137").Value = -12600
$ws.Range("H138").Value = 2892.77
$ws.Range("I138").Value = 1358.1
$ws.Range("J138").Value = 3550.4856
$ws.Range("K138").Value = 4074.3
$ws.Range("L138").Value = 10651.4568
$ws.Range("M138").Value = 1065.7
$ws.Range("N138").Value = -20931.4568

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2542.3606
$ws.Range("I32").Value = 2292.9492
$ws.Range("K32").Value = 2292.9492
$ws.Range("M32").Value = -2005.9492
$ws.Range("H74").Value = 10030.083
$ws.Range("I74").Value = 1895.5
$ws.Range("K74").Value = 1895.5
$ws.Range("M74").Value = -1021.5
$ws.Range("H77").Value = 10030.083
$ws.Range("I77").Value = 1895.5
$ws.Range("K77").Value = 9477.5
$ws.Range("M77").Value = -5109.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 371669.25
$ws.Range("I22").Value = 452.54544
$ws.Range("K22").Value = 452.54544
$ws.Range("M22").Value = -279.54544
$ws.Range("H86").Value = 2164.8333
$ws.Range("I86").Value = 2164.8333
$ws.Range("K86").Value = 2164.8333
$ws.Range("M86").Value = -1041.8333
$ws.Range("H89").Value = 2164.8333
$ws.Range("I89").Value = 2164.8333
$ws.Range("K89").Value = 10824.1665
$ws.Range("M89").Value = -5208.166499999999
$ws.Range("H99").Value = 5192.4287
$ws.Range("I99").Value = 2315.125
$ws.Range("J99").Value = 14399.8
$ws.Range("K99").Value = 2315.125
$ws.Range("L99").Value = 14399.8
$ws.Range("M99").Value = -817.125
$ws.Range("N99").Value = -17395.8
$ws.Range("H105").Value = 4794
$ws.Range("I105").Value = 4607
$ws.Range("J105").Value = 5666.6665
$ws.Range("K105").Value = 4607
$ws.Range("L105").Value = 5666.6665
$ws.Range("M105").Value = -2860
$ws.Range("N105").Value = -9160.666499999999
$ws.Range("H134").Value = 2635.5908
$ws.Range("I134").Value = 2341.2104
$ws.Range("K134").Value = 7023.6312
$ws.Range("M134").Value = -4488.6312

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 352.55554
$ws.Range("I7").Value = 222.375
$ws.Range("K7").Value = 222.375
$ws.Range("M7").Value = -109.375
$ws.Range("H16").Value = 3115.8
$ws.Range("I16").Value = 3115.8
$ws.Range("K16").Value = 3115.8
$ws.Range("M16").Value = -2828.8
$ws.Range("H22").Value = 1742.5714
$ws.Range("I22").Value = 299.75
$ws.Range("J22").Value = 3666.3333
$ws.Range("K22").Value = 299.75
$ws.Range("L22").Value = 3666.3333
$ws.Range("M22").Value = 50.25
$ws.Range("N22").Value = -4366.3333
$ws.Range("H31").Value = 25887.705
$ws.Range("I31").Value = 29641.555
$ws.Range("J31").Value = 8995.375
$ws.Range("K31").Value = 29641.555
$ws.Range("L31").Value = 8995.375
$ws.Range("M31").Value = -29346.555
$ws.Range("N31").Value = -9585.375
$ws.Range("H34").Value = 25887.705
$ws.Range("I34").Value = 29641.555
$ws.Range("J34").Value = 8995.375
$ws.Range("K34").Value = 29641.555
$ws.Range("L34").Value = 8995.375
$ws.Range("M34").Value = -29439.555
$ws.Range("N34").Value = -9399.375
$ws.Range("H58").Value = 2486
$ws.Range("I58").Value = 2699.4285
$ws.Range("K58").Value = 2699.4285
$ws.Range("M58").Value = -2496.4285
$ws.Range("H94").Value = 747.0714
$ws.Range("I94").Value = 703.75
$ws.Range("J94").Value = 804.8333
$ws.Range("K94").Value = 703.75
$ws.Range("L94").Value = 804.8333
$ws.Range("M94").Value = -252.75
$ws.Range("N94").Value = -1706.8333
$ws.Range("H95").Value = 19249.2
$ws.Range("J95").Value = 18561.75
$ws.Range("L95").Value = 18561.75
$ws.Range("N95").Value = -24053.75
$ws.Range("H113").Value = 3115.8
$ws.Range("I113").Value = 3115.8
$ws.Range("K113").Value = 3115.8
$ws.Range("M113").Value = -945.8000000000002
$ws.Range("H132").Value = 3228.2258
$ws.Range("I132").Value = 3040.6897
$ws.Range("J132").Value = 5947.5
$ws.Range("K132").Value = 9122.069100000001
$ws.Range("L132").Value = 17842.5
$ws.Range("M132").Value = -6592.069100000001
$ws.Range("N132").Value = -22902.5
$ws.Range("H136").Value = 2486
$ws.Range("I136").Value = 2699.4285
$ws.Range("K136").Value = 8098.2855
$ws.Range("M136").Value = -5548.2855

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 131.25
$ws.Range("I40").Value = 79.411766
$ws.Range("J40").Value = 257.14285
$ws.Range("K40").Value = 317.647064
$ws.Range("L40").Value = 1028.5714
$ws.Range("M40").Value = -248.647064
$ws.Range("N40").Value = -1166.5714
$ws.Range("H86").Value = 644.5
$ws.Range("I86").Value = 447
$ws.Range("K86").Value = 1341
$ws.Range("M86").Value = -155
$ws.Range("H89").Value = 644.5
$ws.Range("I89").Value = 447
$ws.Range("K89").Value = 4023
$ws.Range("M89").Value = 1905

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 4000
$ws.Range("I59").Value = 4000
$ws.Range("K59").Value = 4000
$ws.Range("M59").Value = -3417
$ws.Range("H113").Value = 5532.364
$ws.Range("I113").Value = 5907.25
$ws.Range("J113").Value = 4532.6665
$ws.Range("K113").Value = 5907.25
$ws.Range("L113").Value = 4532.6665
$ws.Range("M113").Value = -3737.25
$ws.Range("N113").Value = -8872.666499999999
$ws.Range("H124").Value = 25000.75
$ws.Range("J124").Value = 25000.75
$ws.Range("L124").Value = 25000.75
$ws.Range("N124").Value = -34820.75
$ws.Range("H132").Value = 187943.3
$ws.Range("I132").Value = 198918.47
$ws.Range("J132").Value = 1365.3334
$ws.Range("K132").Value = 596755.41
$ws.Range("L132").Value = 4096.0002
$ws.Range("M132").Value = -594225.41
$ws.Range("N132").Value = -9156.0002

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 696875.1
$ws.Range("I93").Value = 743240.1
$ws.Range("K93").Value = 743240.1
$ws.Range("M93").Value = -741992.1
$ws.Range("H132").Value = 3492.5881
$ws.Range("I132").Value = 3148.375
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 9445.125
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -6915.125
$ws.Range("N132").Value = -32060

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 36666
$ws.Range("J63").Value = 36666
$ws.Range("L63").Value = 36666
$ws.Range("N63").Value = -37914
$ws.Range("H66").Value = 36666
$ws.Range("J66").Value = 36666
$ws.Range("L66").Value = 109998
$ws.Range("N66").Value = -116238
$ws.Range("H126").Value = 4828.2666
$ws.Range("I126").Value = 4680.3076
$ws.Range("J126").Value = 5790
$ws.Range("K126").Value = 14040.9228
$ws.Range("L126").Value = 17370
$ws.Range("M126").Value = -11570.9228
$ws.Range("N126").Value = -22310
$ws.Range("H132").Value = 2340.5
$ws.Range("I132").Value = 2177.7856
$ws.Range("J132").Value = 2910
$ws.Range("K132").Value = 6533.3568
$ws.Range("L132").Value = 8730
$ws.Range("M132").Value = -4003.3568
$ws.Range("N132").Value = -13790

